$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet
$ws.Name = "Sheet1"

# Delete "Position Ratio" column (column B) - shifts C->B, D->C, E->D
$ws.Columns.Item(2).Delete()

# Convert all formulas in the used range to static cached values
$used = $ws.UsedRange
$rows = $used.Rows.Count
$cols = $used.Columns.Count
for ($r = 1; $r -le $rows; $r++) {
  for ($c = 1; $c -le $cols; $c++) {
    $cell = $ws.Cells.Item($r, $c)
    $v = $cell.Value()
    if ($null -ne $v) {
      $cell.Value = $v
    }
  }
}

Write-Output "phase1 done"
